$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the numeric-looking values to be stored as text (matching the
# source data, e.g. "Senha"/"CPF"/"CEP"/"Telefone" columns), same as the
# existing rows 2-3. Purely alphabetic/email/address cells do not need
# this - Excel already stores them as text by default.
$ws.Range("B4").NumberFormat = "@"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("G5").NumberFormat = "@"

$ws.Range("A4").Value = "paulo"
$ws.Range("B4").Value = "987584"
$ws.Range("C4").Value = "321321321"
$ws.Range("E4").Value = "92320192"
$ws.Range("F4").Value = "pauloroberto@gmail.com"
$ws.Range("G4").Value = "51999875487"
$ws.Range("H4").Value = "RUA DOS PINHAIS"

$ws.Range("A5").Value = "matheus"
$ws.Range("B5").Value = "9918283182"
$ws.Range("C5").Value = "91283912738216"
$ws.Range("E5").Value = "99827371"
$ws.Range("F5").Value = "matheus@gmail.com"
$ws.Range("G5").Value = "120392391298"
$ws.Range("H5").Value = "rua da topeira"

$ws.Range("A1:H5").Errors(9).Ignore = $true
